$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.985.34'
$ws.Range("E2").Value = '  -5.46%  '
$ws.Range("D3").Value = '2.224.40'
$ws.Range("E3").Value = '  -6.45%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.42%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.583'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.559'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.66'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0823'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -10.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.61'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -10.55%  '
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").Value = '2.561.64'
$ws.Range("E15").Value = '  -6.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.855'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -12.89%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -7.70%  '
$ws.Range("D18").Value = '2.227.57'
$ws.Range("E18").Value = '  -6.45%  '
$ws.Range("D19").Value = '42.899.80'
$ws.Range("E19").Value = '  -5.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.26%  '
$ws.Range("D21").Value = '0.0₃0961'
$ws.Range("E21").Value = '  -9.85%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.50'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -11.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '237.16'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.11'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -12.07%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -11.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0871'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -11.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '34.28'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '155.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.75'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.33%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("B38").Value = 'Stellar'
$ws.Range("C38").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.121'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.92%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.89'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.39'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.68%  '
$ws.Range("E41").Value = '  -11.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -8.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0322'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.58%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.50'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.80%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.791.14'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -13.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.204'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -10.82%  '
$ws.Range("B49").Value = 'ordi'
$ws.Range("C49").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.98'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.68%  '
$ws.Range("B50").Value = 'THORChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.70%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '59.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -15.24%  '
